$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers for new columns F, G, H
$ws.Range("F1").Value = "KNN_Outliers_MAD"
$ws.Range("G1").Value = "SVM_Outliers_MAD"
$ws.Range("H1").Value = "RF_Outliers_MAD"

# Copy header style from existing header (e.g. E1) to the new headers
$ws.Range("E1").Copy()
$ws.Range("F1:H1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A1").Select()

# Outlier flags (MAD-based) per row for KNN (F), SVM (G), RF (H)
$values = @(
    @(2,  $false, $false, $false),
    @(3,  $false, $false, $false),
    @(4,  $false, $false, $false),
    @(5,  $false, $false, $false),
    @(6,  $false, $false, $false),
    @(7,  $false, $false, $false),
    @(8,  $false, $false, $false),
    @(9,  $false, $false, $false),
    @(10, $false, $true,  $false),
    @(11, $false, $false, $false),
    @(12, $false, $false, $false),
    @(13, $false, $false, $false),
    @(14, $false, $false, $false),
    @(15, $false, $false, $false),
    @(16, $false, $false, $false),
    @(17, $false, $false, $false),
    @(18, $false, $false, $false),
    @(19, $false, $false, $false),
    @(20, $true,  $false, $false),
    @(21, $false, $false, $false),
    @(22, $false, $false, $false),
    @(23, $false, $false, $true),
    @(24, $true,  $false, $true),
    @(25, $false, $false, $false)
)

foreach ($row in $values) {
    $r = $row[0]
    $ws.Cells.Item($r, 6).Value = $row[1]
    $ws.Cells.Item($r, 7).Value = $row[2]
    $ws.Cells.Item($r, 8).Value = $row[3]
}
